# Daily attendance processing - reorder "Recorded By" entries in column G
# so that the actual recorder (email / backup) is listed before the
# automated "System"/"system" marker, and "admin@admin.com" is listed
# after other real user entries (e.g. "dnasr281@gmail.com"). The relative
# order of equally-ranked entries is preserved (stable sort).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-TokenRank($tok) {
    if ($tok -eq "System" -or $tok -eq "system") {
        return 2
    } elseif ($tok -eq "admin@admin.com") {
        return 1
    } else {
        return 0
    }
}

function Sort-RecordedByTokens($parts) {
    $n = $parts.Count
    $ranks = @()
    for ($i = 0; $i -lt $n; $i++) {
        $ranks += Get-TokenRank($parts[$i])
    }
    # stable insertion sort (ascending rank): real users first, then
    # admin@admin.com, then System/system last.
    for ($i = 1; $i -lt $n; $i++) {
        $keyRank = $ranks[$i]
        $keyVal = $parts[$i]
        $j = $i - 1
        while ($j -ge 0 -and $ranks[$j] -gt $keyRank) {
            $ranks[$j + 1] = $ranks[$j]
            $parts[$j + 1] = $parts[$j]
            $j = $j - 1
        }
        $ranks[$j + 1] = $keyRank
        $parts[$j + 1] = $keyVal
    }
    return $parts
}

$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1
$firstCol = $ur.Column
$lastCol = $firstCol + $ur.Columns.Count - 1

# Locate the "Recorded By" column dynamically from the header row.
$recordedByCol = 7
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Text
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
    }
}

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByCol)
    $orig = $cell.Text

    if ($orig -ne $null -and $orig.Length -gt 0 -and $orig.Contains(",")) {
        $parts = @($orig -split ", ")
        $sortedParts = Sort-RecordedByTokens($parts)
        $updated = $sortedParts -join ", "

        if ($updated -ne $orig) {
            $cell.Value = $updated
        }
    }
}
